$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: writes $Value into cell $Addr as an exact text string (matching the
# source data's inline-string cells), bypassing Excel's automatic
# number/date auto-detection, and without leaving the temporary text
# NumberFormat applied to the cell once done.
function Set-TextCell([string]$Addr, [string]$Value) {
    $cell = $ws.Range($Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    $cell.Style = "Normal"
}


# --- Rows 27/28, 37/38, 39/40, 47/48: Coin+Link (and Price/Volume) swapped ---
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell "D27" '0.111'
Set-TextCell "E27" '  +0.18%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell "D28" '15.70'
Set-TextCell "E28" '  +0.73%  '

$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell "D37" '0.888'
Set-TextCell "E37" '  +1.95%  '

$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell "D38" '0.932'
Set-TextCell "E38" '  +2.22%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell "D39" '0.0168'
Set-TextCell "E39" '  +0.33%  '

$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell "D40" '0.559'
Set-TextCell "E40" '  -1.17%  '

$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell "D47" '1.78'
Set-TextCell "E47" '  +4.79%  '

$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextCell "D48" '1.786.14'
Set-TextCell "E48" '  +0.69%  '

# --- Remaining rows: refreshed Price (D) and/or Volume(1h) (E) values ---
Set-TextCell "D2" '27.956.90'
Set-TextCell "E2" '  +1.02%  '
Set-TextCell "D3" '1.644.05'
Set-TextCell "E3" '  +1.06%  '
Set-TextCell "E4" '  -0.49%  '
Set-TextCell "D5" '212.37'
Set-TextCell "E5" '  +0.03%  '
Set-TextCell "D6" '0.525'
Set-TextCell "E6" '  +0.42%  '
Set-TextCell "D7" '0.998'
Set-TextCell "E7" '  -0.51%  '
Set-TextCell "D8" '23.56'
Set-TextCell "E8" '  +2.71%  '
Set-TextCell "D9" '0.265'
Set-TextCell "E9" '  +1.22%  '
Set-TextCell "E10" '  +0.27%  '
Set-TextCell "D11" '0.0867'
Set-TextCell "E11" '  -2.53%  '
Set-TextCell "D12" '1.875.47'
Set-TextCell "E12" '  +0.96%  '
Set-TextCell "D13" '1.646.02'
Set-TextCell "E13" '  +1.32%  '
Set-TextCell "E14" '  +0.45%  '
Set-TextCell "E15" '  +1.83%  '
Set-TextCell "D16" '65.64'
Set-TextCell "E16" '  +1.82%  '
Set-TextCell "D17" '27.889.11'
Set-TextCell "E17" '  +0.74%  '
Set-TextCell "D18" '232.35'
Set-TextCell "E18" '  +0.93%  '
Set-TextCell "E19" '  +1.07%  '
Set-TextCell "E20" '  -0.12%  '
Set-TextCell "E21" '  -0.54%  '
Set-TextCell "D22" '10.75'
Set-TextCell "E22" '  +7.73%  '
Set-TextCell "E23" '  +2.14%  '
Set-TextCell "D24" '2.16'
Set-TextCell "E24" '  +3.32%  '
Set-TextCell "D25" '150.80'
Set-TextCell "E25" '  +0.82%  '
Set-TextCell "D26" '6.93'
Set-TextCell "E26" '  +0.60%  '
Set-TextCell "E29" '  -0.49%  '
Set-TextCell "E30" '  +0.29%  '
Set-TextCell "E31" '  -0.02%  '
Set-TextCell "D32" '3.31'
Set-TextCell "E32" '  +0.43%  '
Set-TextCell "D33" '1.468.18'
Set-TextCell "E33" '  +0.35%  '
Set-TextCell "E34" '  +0.45%  '
Set-TextCell "E35" '  +0.82%  '
Set-TextCell "D36" '2.32'
Set-TextCell "E36" '  -1.39%  '
Set-TextCell "D41" '69.22'
Set-TextCell "E41" '  -0.34%  '
Set-TextCell "E42" '  +0.19%  '
Set-TextCell "D43" '0.999'
Set-TextCell "E43" '  -0.45%  '
Set-TextCell "D44" '2.46'
Set-TextCell "E44" '  +0.04%  '
Set-TextCell "E45" '  -0.38%  '
Set-TextCell "D46" '5.39'
Set-TextCell "E46" '  -0.91%  '
Set-TextCell "D49" '87.85'
Set-TextCell "E49" '  +2.09%  '
Set-TextCell "E50" '  +1.48%  '
Set-TextCell "E51" '  -4.31%  '
